# issue #5: stock data output to json file
# Add a "property_category" column (stock sheet) holding the constant
# value "stock" for every row, and tidy up a few company names that had
# a stray embedded space before "股份有限公 司" -> "股份有限公司".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (H) so the
# layout becomes:
#   name, owner, quantity, face_value, currency, total,
#   property_category, date, legislator_name, legislator_id
$ws.Columns.Item(8).Insert()

$lastRow = $ws.UsedRange.Rows.Count

$ws.Range("H1").Value = "property_category"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}

# Fix the three company names that had a stray space before the final
# "司" character.
for ($r = 2; $r -le $lastRow; $r++) {
    $nameCell = $ws.Cells.Item($r, 2)
    $name = $nameCell.Value2
    if ($name -eq "旭晶能源科技股份有限公 司") {
        $nameCell.Value = "旭晶能源科技股份有限公司"
    }
    elseif ($name -eq "珀碩聯.合科技股份有限公 司") {
        $nameCell.Value = "珀碩聯.合科技股份有限公司"
    }
    elseif ($name -eq "景岳生物科技股份有限公 司") {
        $nameCell.Value = "景岳生物科技股份有限公司"
    }
}
